# Edit script: add "Random Forest-100 (superdataset-21.csv)" and
# "...without cons" test blocks (columns N:P and S:U) to sheet1.
# Commit: "superdataset-21 without cons test on mae"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Duplicate the D:F ("train/test (MAE)") and I:K ("...without cons")
#    block layout + styles into N:P and S:U respectively. Row 55 is
#    blank (not part of the used range), so it is copied separately to
#    avoid materialising an empty row that doesn't exist in the source.
# ---------------------------------------------------------------------
$ws.Range("D3:F54").Copy($ws.Range("N3"))
$ws.Range("D56:F57").Copy($ws.Range("N56"))
$ws.Range("I3:K54").Copy($ws.Range("S3"))
$ws.Range("I56:K57").Copy($ws.Range("S56"))

# ---------------------------------------------------------------------
# 2. Fix up the two new group headers.
#    N3 = "Random Forest-100 (superdataset-21.csv without cons)" (has data)
#    S3 = "Random Forest-100 (superdataset-21.csv)" (left empty)
# ---------------------------------------------------------------------
$ws.Range("S3").Value = "Random Forest-100 (superdataset-21.csv)"
$ws.Range("N3").Value = "Random Forest-100 (superdataset-21.csv without cons)"

# ---------------------------------------------------------------------
# 3. Row-index columns (N and S): 1..50 with the same "+1" formula
#    pattern used by D/I.
# ---------------------------------------------------------------------
$ws.Range("N5").Value = 1
$ws.Range("S5").Value = 1
for ($r = 6; $r -le 54; $r++) {
    $ws.Cells.Item($r, 14).Formula = "=N" + ($r - 1) + "+1"
    $ws.Cells.Item($r, 19).Formula = "=S" + ($r - 1) + "+1"
}

# ---------------------------------------------------------------------
# 4. superdataset-21.csv "without cons" train/test MAE values (O, P).
# ---------------------------------------------------------------------
$opValues = @(
    @(53.738101629913707, 149.63628696604599),
    @(53.541086152581833, 155.15138554216861),
    @(54.689241199835642, 143.13598028477551),
    @(55.263889878098887, 147.5578587075575),
    @(54.32085330776605, 150.0395071193866),
    @(53.648561840843733, 151.0301150054764),
    @(55.223594028215309, 144.29275465498361),
    @(56.055545815641693, 144.4272289156626),
    @(54.098357759211062, 148.92030668127049),
    @(53.979420627311328, 144.86430449068999),
    @(54.930138337214082, 144.82167579408539),
    @(54.679001506642933, 149.10549288061341),
    @(54.275036296397751, 148.91654983570649),
    @(55.587652376386792, 137.2727875136911),
    @(54.572720175318452, 150.0026286966046),
    @(55.135577318175592, 142.2148247535597),
    @(56.103305026708661, 137.2438225629792),
    @(54.027246952472261, 152.57814348302301),
    @(54.471531297082599, 144.17922234392111),
    @(54.502900972469533, 148.42997261774369),
    @(54.791131351869602, 146.03740963855421),
    @(54.756736063552943, 145.6442661555312),
    @(54.209411039583607, 145.0289649507119),
    @(55.295178742638001, 136.14734939759029),
    @(54.224766470346523, 151.0955969331873),
    @(54.310672510614992, 141.89998904709751),
    @(53.767460621832619, 152.50776560788611),
    @(54.553750171209423, 141.51168674698789),
    @(54.719635666347081, 141.1243318729463),
    @(54.396527872894133, 141.574282584885),
    @(54.684310368442667, 146.70938116100771),
    @(53.890987535953983, 142.65382803943041),
    @(54.056327900287627, 146.14921686746979),
    @(54.379254896589508, 149.293587075575),
    @(54.591487467470209, 143.4984501642935),
    @(54.169628817970143, 149.89188389923331),
    @(55.281335433502257, 140.73727272727271),
    @(54.485656759348032, 146.27437568455639),
    @(54.15623613203671, 150.5310952902519),
    @(53.809669908231747, 143.680443592552),
    @(53.297362005204761, 150.5702683461117),
    @(55.409305574578823, 139.9791018619934),
    @(53.933847418161882, 144.6483625410734),
    @(53.936604574715787, 149.1498138006572),
    @(53.876227913984387, 154.67331872946329),
    @(54.473872072318862, 146.68242059145669),
    @(56.163641966853866, 130.6130887185104),
    @(55.11714011779209, 143.91588170865279),
    @(54.78606629228873, 142.34866922234389),
    @(55.61784687029175, 144.2873384446878)
)

for ($i = 0; $i -lt $opValues.Length; $i++) {
    $r = 5 + $i
    $ws.Cells.Item($r, 15).Value = $opValues[$i][0]
    $ws.Cells.Item($r, 16).Value = $opValues[$i][1]
}

# ---------------------------------------------------------------------
# 5. superdataset-21.csv (with cons) test was never run -> T/U stay
#    empty (keeps the style copied above, clears the copied values).
# ---------------------------------------------------------------------
$ws.Range("T5:U54").ClearContents()

# ---------------------------------------------------------------------
# 6. Summary rows: avg (56) / SD (57) for both new blocks.
# ---------------------------------------------------------------------
$ws.Range("O56").Formula = "=AVERAGE(O5:O54)"
$ws.Range("P56").Formula = "=AVERAGE(P5:P54)"
$ws.Range("O57").Formula = "=STDEV.S(O5:O54)"
$ws.Range("P57").Formula = "=STDEV.S(P5:P54)"

$ws.Range("T56").Formula = "=AVERAGE(T5:T54)"
$ws.Range("U56").Formula = "=AVERAGE(U5:U54)"
$ws.Range("T57").Formula = "=STDEV.S(T5:T54)"
$ws.Range("U57").Formula = "=STDEV.S(U5:U54)"

# ---------------------------------------------------------------------
# 7. Column widths for the new value columns (O, P, T, U).
# ---------------------------------------------------------------------
$ws.Columns.Item(15).ColumnWidth = 12.307291666666666
$ws.Columns.Item(16).ColumnWidth = 13.022135416666666
$ws.Columns.Item(20).ColumnWidth = 12.736979166666666
$ws.Columns.Item(21).ColumnWidth = 11.451822916666666

# ---------------------------------------------------------------------
# 8. Selection / view: Excel no longer pins topLeftCell, and the last
#    active cell moves to U14.
# ---------------------------------------------------------------------
$ws.Range("U14").Select()
